$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 212, pushing existing rows 212-217 down to 214-219
$ws.Rows("212:213").Insert()

# New row 212 data
$ws.Cells.Item(212,1).Value2 = 9
$ws.Cells.Item(212,2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(212,3).Value2 = "Metropolitana"
$ws.Cells.Item(212,4).Value2 = 44448
$ws.Cells.Item(212,5).Value2 = 13
$ws.Cells.Item(212,6).Value2 = 100112012
$ws.Cells.Item(212,7).Value2 = "Espinaca"
$ws.Cells.Item(212,8).Value2 = "Sin especificar"
$ws.Cells.Item(212,9).Value2 = "Primera"
$ws.Cells.Item(212,10).Value2 = 250
$ws.Cells.Item(212,11).Value2 = 7000
$ws.Cells.Item(212,12).Value2 = 8000
$ws.Cells.Item(212,13).Value2 = 7500
$ws.Cells.Item(212,14).Value2 = "`$/cuna 10 kilos"
$ws.Cells.Item(212,15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(212,16).Value2 = 750
$ws.Cells.Item(212,17).Value2 = 10
$ws.Cells.Item(212,18).Value2 = "Hortaliza"

# New row 213 data
$ws.Cells.Item(213,1).Value2 = 9
$ws.Cells.Item(213,2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(213,3).Value2 = "Metropolitana"
$ws.Cells.Item(213,4).Value2 = 44448
$ws.Cells.Item(213,5).Value2 = 13
$ws.Cells.Item(213,6).Value2 = 100112012
$ws.Cells.Item(213,7).Value2 = "Espinaca"
$ws.Cells.Item(213,8).Value2 = "Sin especificar"
$ws.Cells.Item(213,9).Value2 = "Segunda"
$ws.Cells.Item(213,10).Value2 = 106
$ws.Cells.Item(213,11).Value2 = 5000
$ws.Cells.Item(213,12).Value2 = 6000
$ws.Cells.Item(213,13).Value2 = 5500
$ws.Cells.Item(213,14).Value2 = "`$/cuna 10 kilos"
$ws.Cells.Item(213,15).Value2 = "Provincia de Chacabuco"
$ws.Cells.Item(213,16).Value2 = 550
$ws.Cells.Item(213,17).Value2 = 10
$ws.Cells.Item(213,18).Value2 = "Hortaliza"

# Apply the same date number format style (style index 2 in original file) to column D of new rows
$ws.Cells.Item(212,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
